# Updated cryptos list on Tue Oct 17 18:10:05 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force Excel to store the value as literal text (avoids auto-conversion
    # of number-looking strings like "211.95" into floating point numbers),
    # while keeping the cell's style/format identical to before (no NumberFormat
    # or style index left applied to the cell).
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.489.58"
$ws.Range("E2").Value = "  +0.79%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.570.89"
$ws.Range("E3").Value = "  -0.50%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.29%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "211.95"
$ws.Range("E5").Value = "  -0.63%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.493"
$ws.Range("E6").Value = "  -0.42%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.30%  "

# Row 8 - OKB
Set-TextValue $ws.Range("D8") "46.26"
$ws.Range("E8").Value = "  +6.09%  "

# Row 9 - Solana
Set-TextValue $ws.Range("D9") "24.11"
$ws.Range("E9").Value = "  +2.28%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -1.57%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.55%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.28%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "1.795.72"
$ws.Range("E13").Value = "  -0.44%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.562.24"
$ws.Range("E14").Value = "  -0.98%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("D15") "0.520"
$ws.Range("E15").Value = "  -1.35%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  -2.08%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "28.451.00"
$ws.Range("E17").Value = "  +0.77%  "

# Row 18 - Litecoin
Set-TextValue $ws.Range("D18") "62.03"
$ws.Range("E18").Value = "  -2.84%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "227.08"
$ws.Range("E19").Value = "  -1.81%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "7.34"
$ws.Range("E20").Value = "  -1.74%  "

# Row 21 - ShibaInu: unchanged

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.21%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "3.88"
$ws.Range("E23").Value = "  -5.92%  "

# Row 24 - Avalanche
Set-TextValue $ws.Range("D24") "9.12"
$ws.Range("E24").Value = "  -2.49%  "

# Row 25 - Toncoin
Set-TextValue $ws.Range("D25") "2.05"
$ws.Range("E25").Value = "  +5.06%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "150.87"
$ws.Range("E26").Value = "  -0.71%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "14.97"
$ws.Range("E27").Value = "  -2.06%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  -2.40%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -2.29%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  +0.28%  "

# Row 31 - was Hedera, now PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.11"
$ws.Range("E31").Value = "  -3.37%  "

# Row 32 - was PancakeSwap, now Hedera
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D32") "0.0464"
$ws.Range("E32").Value = "  -2.16%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -0.61%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -0.09%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.389.36"
$ws.Range("E35").Value = "  -2.03%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -3.28%  "

# Row 37 - TrustWalletToken
$ws.Range("E37").Value = "  -2.37%  "

# Row 38 - HuobiToken
$ws.Range("E38").Value = "  +1.54%  "

# Row 39 - MXToken
Set-TextValue $ws.Range("D39") "2.61"
$ws.Range("E39").Value = "  +2.85%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -0.72%  "

# Row 41 - ImmutableX
Set-TextValue $ws.Range("D41") "0.532"
$ws.Range("E41").Value = "  -2.01%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.29%  "

# Row 43 - ARBITRUM
$ws.Range("E43").Value = "  -2.17%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  -1.04%  "

# Row 45 - RenderToken
Set-TextValue $ws.Range("D45") "1.86"
$ws.Range("E45").Value = "  +1.93%  "

# Row 46 - WEMIXToken
Set-TextValue $ws.Range("D46") "0.980"
$ws.Range("E46").Value = "  +0.84%  "

# Row 47 - Aave
Set-TextValue $ws.Range("D47") "62.97"
$ws.Range("E47").Value = "  -1.77%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "1.707.46"
$ws.Range("E48").Value = "  -0.53%  "

# Row 49 - Quant
Set-TextValue $ws.Range("D49") "85.85"
$ws.Range("E49").Value = "  -1.54%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  -3.25%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -1.36%  "
